$d = $word.ActiveDocument

$map = @{
  "99×20=1980" = "81×77=6237"
  "24×90=2160" = "61×64=3904"
  "33×37=1221" = "34×25=850"
  "33×98=3234" = "80×36=2880"
  "93×78=7254" = "65×62=4030"
  "75×22=1650" = "63×58=3654"
  "55×86=4730" = "77×89=6853"
  "27×78=2106" = "49×31=1519"
  "56×87=4872" = "41×35=1435"
  "64×92=5888" = "91×87=7917"
  "49×40=1960" = "98×95=9310"
  "57×16=912"  = "78×59=4602"
  "97×87=8439" = "62×51=3162"
  "46×57=2622" = "76×56=4256"
  "42×19=798"  = "27×66=1782"
  "70×25=1750" = "95×96=9120"
  "28×86=2408" = "72×25=1800"
  "95×62=5890" = "86×22=1892"
  "64×62=3968" = "25×52=1300"
  "83×28=2324" = "74×91=6734"
  "24×41=984"  = "51×95=4845"
  "92×60=5520" = "70×32=2240"
  "75×16=1200" = "33×28=924"
  "90×63=5670" = "58×34=1972"
  "20×89=1780" = "24×97=2328"
}

foreach ($old in $map.Keys) {
  $new = $map[$old]
  $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                           $true, 1, $false, $new, 2)
}
